$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.936.92"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.407.59"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "554.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.06"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "2.399.76"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  +1.13%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.350"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.93"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("E15").Value = "  +5.15%  "
$ws.Range("D16").Value = "2.841.14"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "61.819.40"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "2.402.21"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "322.94"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.17"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "64.99"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +3.36%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.04"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +9.71%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "574.83"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +15.24%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "2.525.28"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.23"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0922"
$ws.Range("E31").Value = "  +5.25%  "
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("E36").Value = "  +0.12%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.62"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +5.99%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.74"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("E39").Value = "  +1.15%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "150.86"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.93%  "
$ws.Range("E41").Value = "  +0.49%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("E43").Value = "  +0.05%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +12.99%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "149.40"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  +1.28%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0536"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.50%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "20.06"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.80%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.586"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.48%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0922"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("E51").Value = "  +2.42%  "
